# Fruta / hortaliza, semanal
# Rotate the weekly price records across rows 3-11 (row 6 unaffected).
# New row N gets the full data previously held by old row Source(N):
#   3 <- 11, 4 <- 10, 5 <- 7, 7 <- 8, 8 <- 4, 9 <- 5, 10 <- 3, 11 <- 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for every row that participates in the rotation,
# across every column that can change (D, H, J, K, L, M, N, O, P).
$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P")
$rowsInvolved = @(3, 4, 5, 7, 8, 9, 10, 11)

$snapshot = @{}
foreach ($r in $rowsInvolved) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: destination row -> source row (values copied from source's "before" snapshot)
$mapping = @{
    3  = 11
    4  = 10
    5  = 7
    7  = 8
    8  = 4
    9  = 5
    10 = 3
    11 = 9
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $srcData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value2 = $srcData[$c]
    }
}
